$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 0.1181316666666667
$ws.Range("I2").Value = 0.1921951611040097
$ws.Range("J2").Value = 0.1921951611040097
$ws.Range("Q2").Value = 0.06460486967444445
$ws.Range("S2").Value = 0.1921951611040097
$ws.Range("T2").Value = 0.1921951611040097

# Row 3 updates
$ws.Range("G3").Value = 0.4965126666666667
$ws.Range("H3").Value = 1.489538
$ws.Range("I3").Value = 0.8078048388959902
$ws.Range("J3").Value = 0.8078048388959903
$ws.Range("Q3").Value = 0.2715371502564444
$ws.Range("R3").Value = 2.443834352308
$ws.Range("S3").Value = 0.8078048388959902
$ws.Range("T3").Value = 0.8078048388959903
